$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.663.52'
$ws.Range("E2").Value = '  +0.96%  '

# Row 3
$ws.Range("D3").Value = '1.641.18'
$ws.Range("E3").Value = '  +0.09%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("E5").Value = '  +0.68%  '

# Row 6
$ws.Range("D6").Value = '''0.527'
$ws.Range("E6").Value = '  -0.80%  '

# Row 7
$ws.Range("E7").Value = '  +0.11%  '

# Row 8
$ws.Range("D8").Value = '''23.22'
$ws.Range("E8").Value = '  +1.00%  '

# Row 9
$ws.Range("E9").Value = '  +0.84%  '

# Row 10
$ws.Range("D10").Value = '''0.0611'
$ws.Range("E10").Value = '  +0.23%  '

# Row 11
$ws.Range("D11").Value = '''0.0896'
$ws.Range("E11").Value = '  +0.36%  '

# Row 12
$ws.Range("D12").Value = '1.874.78'
$ws.Range("E12").Value = '  +0.19%  '

# Row 13
$ws.Range("D13").Value = '1.659.80'
$ws.Range("E13").Value = '  +1.32%  '

# Row 14
$ws.Range("E14").Value = '  +0.55%  '

# Row 15
$ws.Range("D15").Value = '''0.562'
$ws.Range("E15").Value = '  +0.49%  '

# Row 16
$ws.Range("D16").Value = '''64.68'
$ws.Range("E16").Value = '  +0.45%  '

# Row 17
$ws.Range("D17").Value = '27.658.06'
$ws.Range("E17").Value = '  +1.06%  '

# Row 18
$ws.Range("D18").Value = '''231.06'
$ws.Range("E18").Value = '  +0.87%  '

# Row 19
$ws.Range("D19").Value = '''7.71'
$ws.Range("E19").Value = '  +1.94%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0724'
$ws.Range("E20").Value = '  +0.50%  '

# Row 21
$ws.Range("E21").Value = '  +0.06%  '

# Row 22
$ws.Range("E22").Value = '  -0.66%  '

# Row 23
$ws.Range("D23").Value = '''10.05'
$ws.Range("E23").Value = '  +4.79%  '

# Row 24
$ws.Range("E24").Value = '  -3.15%  '

# Row 25
$ws.Range("D25").Value = '''149.94'
$ws.Range("E25").Value = '  +1.96%  '

# Row 26
$ws.Range("D26").Value = '''6.96'
$ws.Range("E26").Value = '  -0.15%  '

# Row 27
$ws.Range("E27").Value = '  -0.84%  '

# Row 28
$ws.Range("D28").Value = '''15.66'
$ws.Range("E28").Value = '  +0.87%  '

# Row 29
$ws.Range("E29").Value = '  +0.08%  '

# Row 30
$ws.Range("D30").Value = '''1.18'
$ws.Range("E30").Value = '  +0.37%  '

# Row 31
$ws.Range("D31").Value = '''0.0487'
$ws.Range("E31").Value = '  +0.70%  '

# Row 32
$ws.Range("E32").Value = '  +0.83%  '

# Row 33
$ws.Range("D33").Value = '1.453.67'
$ws.Range("E33").Value = '  +2.92%  '

# Row 34
$ws.Range("D34").Value = '''3.13'
$ws.Range("E34").Value = '  +0.52%  '

# Row 35
$ws.Range("E35").Value = '  +0.25%  '

# Row 36
$ws.Range("E36").Value = '  -0.86%  '

# Row 37
$ws.Range("D37").Value = '''0.569'
$ws.Range("E37").Value = '  +0.95%  '

# Row 38
$ws.Range("D38").Value = '''0.880'
$ws.Range("E38").Value = '  -0.06%  '

# Row 39
$ws.Range("E39").Value = '  +0.52%  '

# Row 40
$ws.Range("D40").Value = '''0.910'
$ws.Range("E40").Value = '  +14.79%  '

# Row 41
$ws.Range("D41").Value = '''70.59'
$ws.Range("E41").Value = '  +9.44%  '

# Row 42
$ws.Range("E42").Value = '  +0.03%  '

# Row 43
$ws.Range("E43").Value = '  +0.09%  '

# Row 44
$ws.Range("D44").Value = '''5.63'
$ws.Range("E44").Value = '  +2.31%  '

# Row 45
$ws.Range("E45").Value = '  +0.60%  '

# Row 46
$ws.Range("E46").Value = '  +0.76%  '

# Row 47
$ws.Range("D47").Value = '1.784.12'
$ws.Range("E47").Value = '  +0.15%  '

# Row 48
$ws.Range("D48").Value = '''1.71'
$ws.Range("E48").Value = '  +3.32%  '

# Row 49
$ws.Range("D49").Value = '''86.08'
$ws.Range("E49").Value = '  -1.87%  '

# Row 50
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₆0106'
$ws.Range("E50").Value = '  +0.29%  '

# Row 51
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '''0.0990'
$ws.Range("E51").Value = '  +0.42%  '
